$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row re-shuffles (coin order swapped in the source feed) ---
# Cosmos <-> Maker (rows 33/34), Dai <-> Hedera (rows 35/36), TheGraph <-> Stacks (rows 44/45)

# Row 33
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "3.826.25"
$ws.Range("E33").Value = "  -0.21%  "

# Row 34
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.74"
$ws.Range("E34").Value = "  -3.51%  "

# Row 35
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.102"
$ws.Range("E36").Value = "  -3.26%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.322"
$ws.Range("E44").Value = "  -4.23%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.95"
$ws.Range("E45").Value = "  -7.44%  "

# --- Price / Volume(1h) refresh for all other rows ---

$ws.Range("D2").Value = "66.519.13"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.264.47"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.34"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.04"
$ws.Range("E6").Value = "  -7.14%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "3.255.97"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  -5.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.567"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.76"
$ws.Range("E12").Value = "  -5.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "688.93"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "3.789.32"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.21"
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").Value = "66.578.83"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.118"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "3.265.16"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.13"
$ws.Range("E20").Value = "  -4.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.62"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.878"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.78"
$ws.Range("E23").Value = "  -5.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.16"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.66"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.82"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.65"
$ws.Range("E27").Value = "  -4.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.25"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.02"
$ws.Range("E29").Value = "  -4.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.28"
$ws.Range("E30").Value = "  -2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.61"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "568.39"
$ws.Range("E32").Value = "  -6.66%  "
$ws.Range("E37").Value = "  -15.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "54.95"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.127"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.37"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.55"
$ws.Range("E41").Value = "  -4.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "31.19"
$ws.Range("E42").Value = "  -4.41%  "
$ws.Range("D43").Value = "0.0₃0656"
$ws.Range("E43").Value = "  -6.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0398"
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.51"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.35"
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.99"
$ws.Range("E51").Value = "  -1.78%  "
